$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new introductory paragraph right after the blank line
#    that follows the title, and before the requirements table.
# ------------------------------------------------------------------
$introAnchor = $d.Paragraphs(2).Range
$introAnchor.InsertParagraphAfter()

$newPara = $d.Paragraphs(3)
$newPara.Format.SpaceAfter = 0
$newPara.Range.Font.Name = "Arial"
$newPara.Range.Font.Size = 12

$insRange = $newPara.Range
$insRange.Collapse(0)
$insRange.InsertAfter("Please briefly ")
$insRange.Collapse(0)
$insRange.InsertAfter("describe")
$insRange.Collapse(0)
$insRange.InsertAfter(" two requirements for each")
$insRange.Collapse(0)
$insRange.InsertAfter(" service")
$insRange.Collapse(0)
$insRange.InsertAfter(" provided by the web app to the user.  ")
$insRange.Collapse(0)
$insRange.InsertAfter("This description may use natural language, diagrams, or other notations that are understandable to customers. ")
$insRange.Collapse(0)
$insRange.InsertAfter(" ")
$insRange.Collapse(0)
$insRange.InsertAfter("User requirements may be written in broad statements of the required system features or functionality.")

# ------------------------------------------------------------------
# 2) Update the "Website Pages" table-header cell text and fold the
#    old " (Navigation)" run into a single new label.
# ------------------------------------------------------------------
$null = $d.Content.Find.Execute("Website Pages (Navigation)", $true, $false, $false, $false, $false, $true, 1, $false, "Website Pages (Consider site navigation & content presentation)", 2)

# ------------------------------------------------------------------
# 3) Adjust the page margins (top/bottom) on the section.
# ------------------------------------------------------------------
$d.PageSetup.TopMargin = 50.4
$d.PageSetup.BottomMargin = 36
